$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 846.1818
$ws.Cells.Item(2, 9).Value = 390
$ws.Cells.Item(2, 10).Value = 1393.6
$ws.Cells.Item(2, 11).Value = 390
$ws.Cells.Item(2, 12).Value = 1393.6
$ws.Cells.Item(2, 13).Value = -277
$ws.Cells.Item(2, 14).Value = -1619.6
$ws.Cells.Item(4, 8).Value = 1776
$ws.Cells.Item(4, 9).Value = 1561.6
$ws.Cells.Item(4, 10).Value = 1929.1428
$ws.Cells.Item(4, 11).Value = 1561.6
$ws.Cells.Item(4, 12).Value = 1929.1428
$ws.Cells.Item(4, 13).Value = -1447.6
$ws.Cells.Item(4, 14).Value = -2157.1428
$ws.Cells.Item(11, 8).Value = 700.1111
$ws.Cells.Item(11, 9).Value = 700.1111
$ws.Cells.Item(11, 11).Value = 700.1111
$ws.Cells.Item(11, 13).Value = -560.1111
$ws.Cells.Item(12, 8).Value = 753.4
$ws.Cells.Item(12, 9).Value = 445.14285
$ws.Cells.Item(12, 11).Value = 445.14285
$ws.Cells.Item(12, 13).Value = -275.14285
$ws.Cells.Item(17, 8).Value = 54059.895
$ws.Cells.Item(17, 10).Value = 54059.895
$ws.Cells.Item(17, 12).Value = 162179.685
$ws.Cells.Item(17, 14).Value = -162515.685
$ws.Cells.Item(18, 8).Value = 2373.6
$ws.Cells.Item(18, 9).Value = 1296
$ws.Cells.Item(18, 10).Value = 3990
$ws.Cells.Item(18, 11).Value = 1296
$ws.Cells.Item(18, 12).Value = 3990
$ws.Cells.Item(18, 13).Value = -1012
$ws.Cells.Item(18, 14).Value = -4558
$ws.Cells.Item(28, 8).Value = 5506.615
$ws.Cells.Item(28, 10).Value = 7742.5
$ws.Cells.Item(28, 12).Value = 7742.5
$ws.Cells.Item(28, 14).Value = -8712.5
$ws.Cells.Item(32, 8).Value = 2131.7
$ws.Cells.Item(32, 9).Value = 986.5
$ws.Cells.Item(32, 11).Value = 986.5
$ws.Cells.Item(32, 13).Value = -660.5
$ws.Cells.Item(33, 8).Value = 478.66666
$ws.Cells.Item(33, 9).Value = 346.0909
$ws.Cells.Item(33, 10).Value = 843.25
$ws.Cells.Item(33, 11).Value = 346.0909
$ws.Cells.Item(33, 12).Value = 843.25
$ws.Cells.Item(33, 13).Value = -117.0909
$ws.Cells.Item(33, 14).Value = -1301.25
$ws.Cells.Item(81, 8).Value = 99998
$ws.Cells.Item(81, 9).Value = 99998
$ws.Cells.Item(81, 11).Value = 99998
$ws.Cells.Item(81, 13).Value = -99000
$ws.Cells.Item(84, 8).Value = 99998
$ws.Cells.Item(84, 9).Value = 99998
$ws.Cells.Item(84, 11).Value = 299994
$ws.Cells.Item(84, 13).Value = -295002
$ws.Cells.Item(97, 8).Value = 2050
$ws.Cells.Item(97, 10).Value = 2566.6667
$ws.Cells.Item(97, 12).Value = 7700.000100000001
$ws.Cells.Item(97, 14).Value = -8692.000100000001
$ws.Cells.Item(98, 8).Value = 1955.2
$ws.Cells.Item(98, 9).Value = 1533.6666
$ws.Cells.Item(98, 11).Value = 1533.6666
$ws.Cells.Item(98, 13).Value = -35.66660000000002
$ws.Cells.Item(106, 8).Value = 13774.75
$ws.Cells.Item(106, 9).Value = 13774.75
$ws.Cells.Item(106, 11).Value = 13774.75
$ws.Cells.Item(106, 13).Value = -13143.75
$ws.Cells.Item(107, 8).Value = 945
$ws.Cells.Item(107, 9).Value = 892.5
$ws.Cells.Item(107, 11).Value = 892.5
$ws.Cells.Item(107, 13).Value = 1027.5
$ws.Cells.Item(112, 8).Value = 6106.125
$ws.Cells.Item(112, 9).Value = 0
$ws.Cells.Item(112, 10).Value = 6106.125
$ws.Cells.Item(112, 11).Value = 0
$ws.Cells.Item(112, 12).Value = 18318.375
$ws.Cells.Item(112, 13).Value = ""
$ws.Cells.Item(112, 14).Value = -20534.375
$ws.Cells.Item(122, 8).Value = 1955.2
$ws.Cells.Item(122, 9).Value = 1533.6666
$ws.Cells.Item(122, 11).Value = 4600.9998
$ws.Cells.Item(122, 13).Value = -2150.9998
$ws.Cells.Item(132, 8).Value = 2274.4
$ws.Cells.Item(132, 9).Value = 2291.3635
$ws.Cells.Item(132, 11).Value = 6874.0905
$ws.Cells.Item(132, 13).Value = -4344.0905
$ws.Cells.Item(138, 8).Value = 5609.381
$ws.Cells.Item(138, 10).Value = 5993.3125
$ws.Cells.Item(138, 12).Value = 17979.9375
$ws.Cells.Item(138, 14).Value = -28259.9375

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 9892.460999999999
$ws.Cells.Item(32, 9).Value = 9451.055
$ws.Cells.Item(32, 11).Value = 9451.055
$ws.Cells.Item(32, 13).Value = -9164.055
$ws.Cells.Item(45, 8).Value = 4115.3335
$ws.Cells.Item(45, 9).Value = 1375
$ws.Cells.Item(45, 10).Value = 4898.2856
$ws.Cells.Item(45, 11).Value = 1375
$ws.Cells.Item(45, 12).Value = 4898.2856
$ws.Cells.Item(45, 13).Value = -998
$ws.Cells.Item(45, 14).Value = -5652.2856
$ws.Cells.Item(53, 8).Value = 19928
$ws.Cells.Item(53, 9).Value = 9892.5
$ws.Cells.Item(53, 11).Value = 9892.5
$ws.Cells.Item(53, 13).Value = -9210.5
$ws.Cells.Item(61, 8).Value = 8702823
$ws.Cells.Item(61, 9).Value = 15393149
$ws.Cells.Item(61, 10).Value = 5399.7
$ws.Cells.Item(61, 11).Value = 15393149
$ws.Cells.Item(61, 12).Value = 5399.7
$ws.Cells.Item(61, 13).Value = -15392937
$ws.Cells.Item(61, 14).Value = -5823.7
$ws.Cells.Item(97, 8).Value = 1977.0714
$ws.Cells.Item(97, 9).Value = 997.17645
$ws.Cells.Item(97, 10).Value = 3491.4546
$ws.Cells.Item(97, 11).Value = 997.17645
$ws.Cells.Item(97, 12).Value = 3491.4546
$ws.Cells.Item(97, 13).Value = -501.17645
$ws.Cells.Item(97, 14).Value = -4483.4546
$ws.Cells.Item(102, 8).Value = 2932.6316
$ws.Cells.Item(102, 9).Value = 2295.2942
$ws.Cells.Item(102, 10).Value = 8350
$ws.Cells.Item(102, 11).Value = 2295.2942
$ws.Cells.Item(102, 12).Value = 8350
$ws.Cells.Item(102, 13).Value = -673.2941999999998
$ws.Cells.Item(102, 14).Value = -11594
$ws.Cells.Item(122, 8).Value = 6538.6665
$ws.Cells.Item(122, 9).Value = 5164.6665
$ws.Cells.Item(122, 10).Value = 8599.666999999999
$ws.Cells.Item(122, 11).Value = 15493.9995
$ws.Cells.Item(122, 12).Value = 25799.001
$ws.Cells.Item(122, 13).Value = -13043.9995
$ws.Cells.Item(122, 14).Value = -30699.001
$ws.Cells.Item(132, 8).Value = 7049.125
$ws.Cells.Item(132, 9).Value = 7049.125
$ws.Cells.Item(132, 11).Value = 21147.375
$ws.Cells.Item(132, 13).Value = -18617.375
$ws.Cells.Item(136, 8).Value = 8702823
$ws.Cells.Item(136, 9).Value = 15393149
$ws.Cells.Item(136, 10).Value = 5399.7
$ws.Cells.Item(136, 11).Value = 46179447
$ws.Cells.Item(136, 12).Value = 16199.1
$ws.Cells.Item(136, 13).Value = -46176897
$ws.Cells.Item(136, 14).Value = -21299.1

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 3780.1177
$ws.Cells.Item(86, 9).Value = 2370.0908
$ws.Cells.Item(86, 10).Value = 6365.1665
$ws.Cells.Item(86, 11).Value = 2370.0908
$ws.Cells.Item(86, 12).Value = 6365.1665
$ws.Cells.Item(86, 13).Value = -1247.0908
$ws.Cells.Item(86, 14).Value = -8611.166499999999
$ws.Cells.Item(89, 8).Value = 3780.1177
$ws.Cells.Item(89, 9).Value = 2370.0908
$ws.Cells.Item(89, 10).Value = 6365.1665
$ws.Cells.Item(89, 11).Value = 11850.454
$ws.Cells.Item(89, 12).Value = 31825.8325
$ws.Cells.Item(89, 13).Value = -6234.454
$ws.Cells.Item(89, 14).Value = -43057.8325
$ws.Cells.Item(94, 8).Value = 1991.5667
$ws.Cells.Item(94, 9).Value = 2093.739
$ws.Cells.Item(94, 11).Value = 2093.739
$ws.Cells.Item(94, 13).Value = -1642.739
$ws.Cells.Item(99, 8).Value = 2784.2856
$ws.Cells.Item(99, 9).Value = 1250
$ws.Cells.Item(99, 11).Value = 1250
$ws.Cells.Item(99, 13).Value = 248
$ws.Cells.Item(105, 8).Value = 922122.4
$ws.Cells.Item(105, 9).Value = 1428656.1
$ws.Cells.Item(105, 11).Value = 1428656.1
$ws.Cells.Item(105, 13).Value = -1426909.1
$ws.Cells.Item(107, 8).Value = 7856.5713
$ws.Cells.Item(107, 9).Value = 7856.5713
$ws.Cells.Item(107, 11).Value = 7856.5713
$ws.Cells.Item(107, 13).Value = -5936.5713
$ws.Cells.Item(134, 8).Value = 3662
$ws.Cells.Item(134, 9).Value = 2550
$ws.Cells.Item(134, 11).Value = 7650
$ws.Cells.Item(134, 13).Value = -5115
$ws.Cells.Item(139, 8).Value = 124374
$ws.Cells.Item(139, 10).Value = 124374
$ws.Cells.Item(139, 12).Value = 124374
$ws.Cells.Item(139, 14).Value = -134654

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 886.5217
$ws.Cells.Item(16, 9).Value = 596
$ws.Cells.Item(16, 10).Value = 1709.6666
$ws.Cells.Item(16, 11).Value = 596
$ws.Cells.Item(16, 12).Value = 1709.6666
$ws.Cells.Item(16, 13).Value = -309
$ws.Cells.Item(16, 14).Value = -2283.6666
$ws.Cells.Item(31, 8).Value = 18870956
$ws.Cells.Item(31, 9).Value = 23812094
$ws.Cells.Item(31, 10).Value = 4787.364
$ws.Cells.Item(31, 11).Value = 23812094
$ws.Cells.Item(31, 12).Value = 4787.364
$ws.Cells.Item(31, 13).Value = -23811799
$ws.Cells.Item(31, 14).Value = -5377.364
$ws.Cells.Item(32, 8).Value = 8499.5
$ws.Cells.Item(32, 9).Value = 8499.5
$ws.Cells.Item(32, 11).Value = 8499.5
$ws.Cells.Item(32, 13).Value = -8183.5
$ws.Cells.Item(34, 8).Value = 18870956
$ws.Cells.Item(34, 9).Value = 23812094
$ws.Cells.Item(34, 10).Value = 4787.364
$ws.Cells.Item(34, 11).Value = 23812094
$ws.Cells.Item(34, 12).Value = 4787.364
$ws.Cells.Item(34, 13).Value = -23811892
$ws.Cells.Item(34, 14).Value = -5191.364
$ws.Cells.Item(94, 8).Value = 1899
$ws.Cells.Item(94, 9).Value = 1198.6666
$ws.Cells.Item(94, 10).Value = 2599.3333
$ws.Cells.Item(94, 11).Value = 1198.6666
$ws.Cells.Item(94, 12).Value = 2599.3333
$ws.Cells.Item(94, 13).Value = -747.6666
$ws.Cells.Item(94, 14).Value = -3501.3333
$ws.Cells.Item(99, 8).Value = 21305
$ws.Cells.Item(99, 9).Value = 12055.857
$ws.Cells.Item(99, 11).Value = 12055.857
$ws.Cells.Item(99, 13).Value = -10557.857
$ws.Cells.Item(105, 8).Value = 5183.5884
$ws.Cells.Item(105, 9).Value = 1009.46155
$ws.Cells.Item(105, 11).Value = 1009.46155
$ws.Cells.Item(105, 13).Value = 737.53845
$ws.Cells.Item(107, 8).Value = 1752.6
$ws.Cells.Item(107, 9).Value = 534.3333
$ws.Cells.Item(107, 11).Value = 534.3333
$ws.Cells.Item(107, 13).Value = 1385.6667
$ws.Cells.Item(113, 8).Value = 886.5217
$ws.Cells.Item(113, 9).Value = 596
$ws.Cells.Item(113, 10).Value = 1709.6666
$ws.Cells.Item(113, 11).Value = 596
$ws.Cells.Item(113, 12).Value = 1709.6666
$ws.Cells.Item(113, 13).Value = 1574
$ws.Cells.Item(113, 14).Value = -6049.6666
$ws.Cells.Item(126, 8).Value = 21305
$ws.Cells.Item(126, 9).Value = 12055.857
$ws.Cells.Item(126, 11).Value = 36167.571
$ws.Cells.Item(126, 13).Value = -33697.571
$ws.Cells.Item(132, 8).Value = 532.2
$ws.Cells.Item(132, 9).Value = 290.25
$ws.Cells.Item(132, 10).Value = 1500
$ws.Cells.Item(132, 11).Value = 870.75
$ws.Cells.Item(132, 12).Value = 4500
$ws.Cells.Item(132, 13).Value = 1659.25
$ws.Cells.Item(132, 14).Value = -9560
$ws.Cells.Item(134, 8).Value = 2594.2856
$ws.Cells.Item(134, 9).Value = 2509.4736
$ws.Cells.Item(134, 11).Value = 7528.4208
$ws.Cells.Item(134, 13).Value = -4993.4208

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 4960.5713
$ws.Cells.Item(12, 9).Value = 4312
$ws.Cells.Item(12, 10).Value = 5220
$ws.Cells.Item(12, 11).Value = 12936
$ws.Cells.Item(12, 12).Value = 15660
$ws.Cells.Item(12, 13).Value = -12763
$ws.Cells.Item(12, 14).Value = -16006
$ws.Cells.Item(23, 8).Value = 22223110
$ws.Cells.Item(23, 9).Value = 210.5
$ws.Cells.Item(23, 11).Value = 631.5
$ws.Cells.Item(23, 13).Value = -396.5
$ws.Cells.Item(50, 8).Value = 882.3333
$ws.Cells.Item(50, 9).Value = 575.75
$ws.Cells.Item(50, 10).Value = 1495.5
$ws.Cells.Item(50, 11).Value = 1727.25
$ws.Cells.Item(50, 12).Value = 4486.5
$ws.Cells.Item(50, 13).Value = -1246.25
$ws.Cells.Item(50, 14).Value = -5448.5
$ws.Cells.Item(53, 8).Value = 882.3333
$ws.Cells.Item(53, 9).Value = 575.75
$ws.Cells.Item(53, 10).Value = 1495.5
$ws.Cells.Item(53, 11).Value = 1727.25
$ws.Cells.Item(53, 12).Value = 4486.5
$ws.Cells.Item(53, 13).Value = -1246.25
$ws.Cells.Item(53, 14).Value = -5448.5
$ws.Cells.Item(63, 8).Value = 27146
$ws.Cells.Item(63, 9).Value = 20011.5
$ws.Cells.Item(63, 11).Value = 60034.5
$ws.Cells.Item(63, 13).Value = -59285.5
$ws.Cells.Item(66, 8).Value = 27146
$ws.Cells.Item(66, 9).Value = 20011.5
$ws.Cells.Item(66, 11).Value = 180103.5
$ws.Cells.Item(66, 13).Value = -176359.5
$ws.Cells.Item(101, 8).Value = 27993
$ws.Cells.Item(101, 10).Value = 27993
$ws.Cells.Item(101, 12).Value = 83979
$ws.Cells.Item(101, 14).Value = -88847
$ws.Cells.Item(110, 8).Value = 18000
$ws.Cells.Item(110, 9).Value = 18000
$ws.Cells.Item(110, 11).Value = 54000
$ws.Cells.Item(110, 13).Value = -49910

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 5600.1
$ws.Cells.Item(80, 9).Value = 3927.8572
$ws.Cells.Item(80, 10).Value = 9502
$ws.Cells.Item(80, 11).Value = 3927.8572
$ws.Cells.Item(80, 12).Value = 9502
$ws.Cells.Item(80, 13).Value = -2929.8572
$ws.Cells.Item(80, 14).Value = -11498
$ws.Cells.Item(83, 8).Value = 5600.1
$ws.Cells.Item(83, 9).Value = 3927.8572
$ws.Cells.Item(83, 10).Value = 9502
$ws.Cells.Item(83, 11).Value = 19639.286
$ws.Cells.Item(83, 12).Value = 47510
$ws.Cells.Item(83, 13).Value = -14647.286
$ws.Cells.Item(83, 14).Value = -57494
$ws.Cells.Item(97, 8).Value = 1323.56
$ws.Cells.Item(97, 9).Value = 1457.2778
$ws.Cells.Item(97, 10).Value = 979.7143
$ws.Cells.Item(97, 11).Value = 1457.2778
$ws.Cells.Item(97, 12).Value = 979.7143
$ws.Cells.Item(97, 13).Value = -961.2778000000001
$ws.Cells.Item(97, 14).Value = -1971.7143
$ws.Cells.Item(102, 8).Value = 4991.8335
$ws.Cells.Item(102, 9).Value = 4987.75
$ws.Cells.Item(102, 11).Value = 4987.75
$ws.Cells.Item(102, 13).Value = -3365.75
$ws.Cells.Item(122, 8).Value = 4221.931
$ws.Cells.Item(122, 9).Value = 4792.8423
$ws.Cells.Item(122, 11).Value = 14378.5269
$ws.Cells.Item(122, 13).Value = -11928.5269
$ws.Cells.Item(126, 8).Value = 12081207
$ws.Cells.Item(126, 10).Value = 8049.1665
$ws.Cells.Item(126, 12).Value = 24147.4995
$ws.Cells.Item(126, 14).Value = -29087.4995
$ws.Cells.Item(132, 8).Value = 11666.5
$ws.Cells.Item(132, 9).Value = 11997
$ws.Cells.Item(132, 10).Value = 10014
$ws.Cells.Item(132, 11).Value = 35991
$ws.Cells.Item(132, 12).Value = 30042
$ws.Cells.Item(132, 13).Value = -33461
$ws.Cells.Item(132, 14).Value = -35102

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1825.7693
$ws.Cells.Item(16, 9).Value = 1825.7693
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 11).Value = 1825.7693
$ws.Cells.Item(16, 12).Value = 0
$ws.Cells.Item(16, 13).Value = -1655.7693
$ws.Cells.Item(16, 14).Value = ""
$ws.Cells.Item(18, 8).Value = 0
$ws.Cells.Item(18, 10).Value = 0
$ws.Cells.Item(18, 12).Value = 0
$ws.Cells.Item(18, 14).Value = ""
$ws.Cells.Item(20, 8).Value = 7000
$ws.Cells.Item(20, 9).Value = 7000
$ws.Cells.Item(20, 11).Value = 7000
$ws.Cells.Item(20, 13).Value = -6774
$ws.Cells.Item(22, 8).Value = 18858464
$ws.Cells.Item(22, 10).Value = 7900
$ws.Cells.Item(22, 12).Value = 7900
$ws.Cells.Item(22, 14).Value = -8490
$ws.Cells.Item(27, 8).Value = 18858464
$ws.Cells.Item(27, 10).Value = 7900
$ws.Cells.Item(27, 12).Value = 7900
$ws.Cells.Item(27, 14).Value = -8114
$ws.Cells.Item(46, 8).Value = 1225.4445
$ws.Cells.Item(46, 9).Value = 1041
$ws.Cells.Item(46, 11).Value = 1041
$ws.Cells.Item(46, 13).Value = -853
$ws.Cells.Item(68, 8).Value = 5212301.5
$ws.Cells.Item(68, 10).Value = 11999.5
$ws.Cells.Item(68, 12).Value = 11999.5
$ws.Cells.Item(68, 14).Value = -13497.5
$ws.Cells.Item(71, 8).Value = 5212301.5
$ws.Cells.Item(71, 10).Value = 11999.5
$ws.Cells.Item(71, 12).Value = 59997.5
$ws.Cells.Item(71, 14).Value = -67485.5
$ws.Cells.Item(82, 8).Value = 4498.3335
$ws.Cells.Item(82, 9).Value = 1614.6666
$ws.Cells.Item(82, 10).Value = 10265.667
$ws.Cells.Item(82, 11).Value = 1614.6666
$ws.Cells.Item(82, 12).Value = 10265.667
$ws.Cells.Item(82, 13).Value = -1253.6666
$ws.Cells.Item(82, 14).Value = -10987.667
$ws.Cells.Item(85, 8).Value = 4498.3335
$ws.Cells.Item(85, 9).Value = 1614.6666
$ws.Cells.Item(85, 10).Value = 10265.667
$ws.Cells.Item(85, 11).Value = 1614.6666
$ws.Cells.Item(85, 12).Value = 10265.667
$ws.Cells.Item(85, 13).Value = -366.6666
$ws.Cells.Item(85, 14).Value = -12761.667
$ws.Cells.Item(93, 8).Value = 2646512.8
$ws.Cells.Item(93, 9).Value = 748.73334
$ws.Cells.Item(93, 10).Value = 9260922
$ws.Cells.Item(93, 11).Value = 748.73334
$ws.Cells.Item(93, 12).Value = 9260922
$ws.Cells.Item(93, 13).Value = 499.26666
$ws.Cells.Item(93, 14).Value = -9263418
$ws.Cells.Item(100, 8).Value = 20856454
$ws.Cells.Item(100, 9).Value = 2459.4
$ws.Cells.Item(100, 10).Value = 35752164
$ws.Cells.Item(100, 11).Value = 2459.4
$ws.Cells.Item(100, 12).Value = 35752164
$ws.Cells.Item(100, 13).Value = -1918.4
$ws.Cells.Item(100, 14).Value = -35753246
$ws.Cells.Item(112, 8).Value = 130899
$ws.Cells.Item(112, 10).Value = 130899
$ws.Cells.Item(112, 12).Value = 130899
$ws.Cells.Item(112, 14).Value = -133853
$ws.Cells.Item(132, 8).Value = 4513.6875
$ws.Cells.Item(132, 9).Value = 3270.2222
$ws.Cells.Item(132, 11).Value = 9810.6666
$ws.Cells.Item(132, 13).Value = -7280.6666
$ws.Cells.Item(136, 8).Value = 4466.1816
$ws.Cells.Item(136, 9).Value = 4375.857
$ws.Cells.Item(136, 11).Value = 13127.571
$ws.Cells.Item(136, 13).Value = -10577.571

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(16, 8).Value = 94715
$ws.Cells.Item(16, 10).Value = 94715
$ws.Cells.Item(16, 12).Value = 94715
$ws.Cells.Item(16, 14).Value = -95299
$ws.Cells.Item(62, 8).Value = 13837
$ws.Cells.Item(62, 9).Value = 7850
$ws.Cells.Item(62, 10).Value = 15832.667
$ws.Cells.Item(62, 11).Value = 7850
$ws.Cells.Item(62, 12).Value = 15832.667
$ws.Cells.Item(62, 13).Value = -7226
$ws.Cells.Item(62, 14).Value = -17080.667
$ws.Cells.Item(65, 8).Value = 13837
$ws.Cells.Item(65, 9).Value = 7850
$ws.Cells.Item(65, 10).Value = 15832.667
$ws.Cells.Item(65, 11).Value = 39250
$ws.Cells.Item(65, 12).Value = 79163.33499999999
$ws.Cells.Item(65, 13).Value = -36130
$ws.Cells.Item(65, 14).Value = -85403.33499999999
$ws.Cells.Item(75, 8).Value = 0
$ws.Cells.Item(75, 10).Value = 0
$ws.Cells.Item(75, 12).Value = 0
$ws.Cells.Item(75, 14).Value = ""
$ws.Cells.Item(78, 8).Value = 0
$ws.Cells.Item(78, 10).Value = 0
$ws.Cells.Item(78, 12).Value = 0
$ws.Cells.Item(78, 14).Value = ""
$ws.Cells.Item(96, 8).Value = 1173.6
$ws.Cells.Item(96, 9).Value = 974.6667
$ws.Cells.Item(96, 10).Value = 1472
$ws.Cells.Item(96, 11).Value = 974.6667
$ws.Cells.Item(96, 12).Value = 1472
$ws.Cells.Item(96, 13).Value = 398.3333
$ws.Cells.Item(96, 14).Value = -4218
$ws.Cells.Item(100, 8).Value = 1542.7059
$ws.Cells.Item(100, 9).Value = 946.8461
$ws.Cells.Item(100, 10).Value = 3479.25
$ws.Cells.Item(100, 11).Value = 1893.6922
$ws.Cells.Item(100, 12).Value = 6958.5
$ws.Cells.Item(100, 13).Value = -1352.6922
$ws.Cells.Item(100, 14).Value = -8040.5
$ws.Cells.Item(107, 8).Value = 5674.9565
$ws.Cells.Item(107, 9).Value = 2944.5833
$ws.Cells.Item(107, 11).Value = 8833.749899999999
$ws.Cells.Item(107, 13).Value = -6913.749899999999
$ws.Cells.Item(122, 8).Value = 3322.25
$ws.Cells.Item(122, 10).Value = 5466.6665
$ws.Cells.Item(122, 12).Value = 16399.9995
$ws.Cells.Item(122, 14).Value = -21299.9995
$ws.Cells.Item(126, 8).Value = 6538.769
$ws.Cells.Item(126, 9).Value = 5191.273
$ws.Cells.Item(126, 10).Value = 13950
$ws.Cells.Item(126, 11).Value = 15573.819
$ws.Cells.Item(126, 12).Value = 41850
$ws.Cells.Item(126, 13).Value = -13103.819
$ws.Cells.Item(126, 14).Value = -46790
$ws.Cells.Item(132, 8).Value = 4415.7144
$ws.Cells.Item(132, 9).Value = 2227
$ws.Cells.Item(132, 10).Value = 7334
$ws.Cells.Item(132, 11).Value = 6681
$ws.Cells.Item(132, 12).Value = 22002
$ws.Cells.Item(132, 13).Value = -4151
$ws.Cells.Item(132, 14).Value = -27062
$ws.Cells.Item(136, 8).Value = 4535.1665
$ws.Cells.Item(136, 9).Value = 5357.143
$ws.Cells.Item(136, 10).Value = 3384.4
$ws.Cells.Item(136, 11).Value = 16071.429
$ws.Cells.Item(136, 12).Value = 10153.2
$ws.Cells.Item(136, 13).Value = -13521.429
$ws.Cells.Item(136, 14).Value = -15253.2
